{"js": "const body = context.document.body;\nbody.insertParagraph(\"test\", \"End\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$ccs = $d.ContentControls\n$cc = $ccs.Item(1)\nWrite-Output \"Type: $($cc.Type)\"\nWrite-Output \"Title: $($cc.Title)\"\nWrite-Output \"Tag: $($cc.Tag)\"\nWrite-Output \"ID: $($cc.ID)\"\nWrite-Output \"BuildingBlockType: $($cc.BuildingBlockType)\"\nWrite-Output \"BuildingBlockCategory: $($cc.BuildingBlockCategory)\"\n"}
